$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 14
$ws.Cells.Item($row, 1).Value = 42620.889479166668
$ws.Cells.Item($row, 2).Value = 18
$ws.Cells.Item($row, 3).Value = 62
$ws.Cells.Item($row, 4).Value = 37
$ws.Cells.Item($row, 5).Value = 62
$ws.Cells.Item($row, 6).Value = 18
$ws.Cells.Item($row, 7).Value = 27764
$ws.Cells.Item($row, 8).Value = 12928
$ws.Cells.Item($row, 9).Value = 2174
$ws.Cells.Item($row, 10).Value = 272
$ws.Cells.Item($row, 11).Value = 162
$ws.Cells.Item($row, 12).Value = 22
$ws.Cells.Item($row, 13).Value = 5
$ws.Cells.Item($row, 14).Value = "Noun"
